$d = $word.ActiveDocument

# Remove the " – V1" suffix from the title paragraph
$d.Content.Find.Execute(" – V1", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Color the "Requisitos do sistema:" paragraph (including its paragraph mark) red
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Color = 255
